$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-14
# from serial 45170 (2023-09-01) to serial 45174 (2023-09-05)
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
